$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("KHADIJA PAPA", "K2546546", "345646131418647463156165", "RABAT 562", "BMCE", "Direction régionale", "800/RABAT AGDAL", "mensuelle", 24000, 1800, 22200),
    @("JAJA GAGA", "B3541456", "354657464131354681448831", "CASA B2", "BP", "Logement de fonction", "800/LF/RABAT AGDAL", "trimestrielle", 40000, 8400, 37200),
    @("MANAL LALA", "G364861", "134165465131864864135418", "UHIU528", "BP", "Point de vente", "622/CASA MEDINA", "annuelle", 150000, 270000, 127500)
)

$startRow = 6
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 8; $col++) {
        # Column C (index 3) holds long numeric account numbers that must
        # stay as text (mirrors the "number stored as text" cells already
        # present in the sheet), so force a text format before assigning.
        if ($col -eq 3) {
            $ws.Cells.Item($row, $col).NumberFormat = "@"
        }
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
    for ($col = 9; $col -le 11; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}

# Extend the "numbers stored as text" ignored-error range (used for the
# CIN/account-number columns) to cover the newly added rows, same as
# Excel does when you click "Ignore Error" across the whole used range.
$endRow = $startRow + $data.Length - 1
$ws.Range("A1:K$endRow").Errors.Item(3).Ignore = $true

